$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()
$ws.Columns.Item(14).EntireColumn.Insert()
$ws.Columns.Item(14).ColumnWidth = 10.33
